# Fixing completion Order bug and tweak order timing
# The "currentshop" sheet's item list is reordered/renumbered: the old
# rows 22-37 (Coffee Stand ... Steak Alarm) become the new rows 2-17,
# with a handful of Slot values corrected (rows 14-17), and the old
# rows 2-21 / 18-37 tail is dropped entirely (dimension shrinks to D17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currentshop")

$items = @(
    @("Coffee Stand",     65,   8,  7),
    @("Newspaper Stand",  80,   4,  7),
    @("Hypno Clock",      80,  15,  7),
    @("Gumball Machine", 120,   6,  7),
    @("Extra Burner",    150,  13,  7),
    @("Ceiling Fan",     150,   2,  7),
    @("TV",               150,  3,  7),
    @("Extra Burner 2",  200,  14,  7),
    @("Arcade Cabinet",  400,   5,  7),
    @("Jukebox",         500,   7,  7),
    @("Royal Crown",    1000,  16,  7),
    @("Doorbell",         30,   1,  8),
    @("Beef Alarm",       90,   9,  8),
    @("Chicken Alarm",    90,  10,  8),
    @("Pork Alarm",       90,  11,  8),
    @("Steak Alarm",      90,  12,  8)
)

for ($i = 0; $i -lt $items.Length; $i++) {
    $row = 2 + $i
    $entry = $items[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}

# Drop the now-duplicated tail (old rows 18-37) so the sheet ends at row 17.
$ws.Rows("18:37").Delete()
